$d = $word.ActiveDocument

# 1. Update the letter date: September 19, 2025 -> September 21, 2025
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing-address paragraph "929 Story Road, San Jose CA 95122"
#    into two paragraphs: "929 Story Road" and "San Jose, CA 95122".
#    (Only the first occurrence, in the letterhead -- not the one inside the
#    "PROPERTY ADDRESS" table further down, so search from the very start of
#    the document and take whichever paragraph the match falls in.)
$rng = $d.Content
$rng.Find.Execute("929 Story Road, San Jose CA 95122") | Out-Null
$addrPara = $rng.Paragraphs(1)
$addrPara.Range.InsertParagraphAfter()
$newPara = $addrPara.Next()
$newPara.Range.Text = "San Jose, CA 95122"
$addrPara.Range.Text = "929 Story Road"

# 3. Remove the empty "No Spacing" paragraph that follows
#    "...Board of Directors".
$rng2 = $d.Content
$rng2.Find.Execute("Board of Directors") | Out-Null
$bodPara = $rng2.Paragraphs(1)
$emptyPara = $bodPara.Next()
$emptyPara.Range.Delete()
